# Update the F column (dSF) values on Sheet1 to reflect repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    2  = -1
    3  = -2
    5  = -1
    6  = 2
    7  = -2
    8  = -4
    10 = -1
    11 = 1
    12 = 3
    13 = -1
    14 = 1
    15 = 2
    16 = -3
    17 = -2
    18 = 6
    20 = 3
    21 = 1
    23 = -1
    24 = 1
    25 = -3
    26 = 3
    27 = -4
    28 = 1
    30 = 1
    31 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
